# Week 13 logging update
# - Rushing sheet: insert a new row for A.St. Brown (between J.Cabinda and K.Raymond),
#   shifting the trailing J.Jefferson row down, and update several stat values.
# - Receiving sheet: update stat values for several players (no row/player changes).

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Sheet 1: Rushing
# ----------------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# Make room for the new "A.St. Brown" row by pushing the last row (old row 8,
# the second J.Jefferson entry) down to row 9. Copy formatting from the row
# above (row 7, same style as every data row) onto the new row 9 first so the
# shifted row keeps the bordered/bold index-cell styling, then fill in values.
$rushing.Range("A7:F7").Copy()
$rushing.Range("A9:F9").PasteSpecial(-4122)

# New row 9 = old row 8 data (J.Jefferson), with updated week-13 totals.
$rushing.Range("A9").Value = 7
$rushing.Range("B9").Value = "J.Jefferson"
$rushing.Range("C9").Value = 3
$rushing.Range("D9").Value = 3
$rushing.Range("E9").Value = 1
$rushing.Range("F9").Value = 2

# Row 7 becomes the new "A.St. Brown" entry.
$rushing.Range("A7").Value = 5
$rushing.Range("B7").Value = "A.St. Brown"
$rushing.Range("C7").Value = 1
$rushing.Range("D7").Value = 0
$rushing.Range("E7").Value = 0
$rushing.Range("F7").Value = 0

# Row 8 is the K.Raymond row (formerly row 7); index shifts down but stats
# are unchanged.
$rushing.Range("A8").Value = 6
$rushing.Range("B8").Value = "K.Raymond"
$rushing.Range("C8").Value = 0
$rushing.Range("D8").Value = 2
$rushing.Range("E8").Value = 0
$rushing.Range("F8").Value = 0

# Updated week-13 totals for the existing rows above the insert point.
$rushing.Range("A2").Value = 0
$rushing.Range("C2").Value = 8
$rushing.Range("E2").Value = 4

$rushing.Range("A3").Value = 1
$rushing.Range("C3").Value = 64
$rushing.Range("D3").Value = 32
$rushing.Range("E3").Value = 15
$rushing.Range("F3").Value = 13

$rushing.Range("A4").Value = 2

$rushing.Range("A5").Value = 3
$rushing.Range("C5").Value = 4
$rushing.Range("D5").Value = 2

$rushing.Range("A6").Value = 4

# ----------------------------------------------------------------------------
# Sheet 2: Receiving
# ----------------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

$receiving.Range("A2").Value = 0
$receiving.Range("C2").Value = 24
$receiving.Range("D2").Value = 22

$receiving.Range("A3").Value = 1

$receiving.Range("A4").Value = 2
$receiving.Range("C4").Value = 7
$receiving.Range("D4").Value = 5

$receiving.Range("A5").Value = 3
$receiving.Range("C5").Value = 2

$receiving.Range("A6").Value = 4

$receiving.Range("A7").Value = 5
$receiving.Range("C7").Value = 26
$receiving.Range("D7").Value = 23
$receiving.Range("E7").Value = 5
$receiving.Range("F7").Value = 2
$receiving.Range("G7").Value = 3
$receiving.Range("H7").Value = 2

$receiving.Range("A8").Value = 6
$receiving.Range("C8").Value = 37
$receiving.Range("D8").Value = 27
$receiving.Range("E8").Value = 15
$receiving.Range("F8").Value = 7

$receiving.Range("A9").Value = 7

$receiving.Range("A10").Value = 8
$receiving.Range("C10").Value = 11
$receiving.Range("D10").Value = 7

$receiving.Range("A11").Value = 9

$receiving.Range("A12").Value = 10

$receiving.Range("A13").Value = 11
$receiving.Range("C13").Value = 5

$receiving.Range("A14").Value = 12
$receiving.Range("C14").Value = 5
$receiving.Range("D14").Value = 3
$receiving.Range("E14").Value = 1
$receiving.Range("F14").Value = 1

$receiving.Range("A15").Value = 13
$receiving.Range("C15").Value = 6
$receiving.Range("D15").Value = 4
$receiving.Range("E15").Value = 6
$receiving.Range("F15").Value = 3
$receiving.Range("G15").Value = 1

$receiving.Range("A16").Value = 14
$receiving.Range("C16").Value = 71
$receiving.Range("D16").Value = 53
$receiving.Range("E16").Value = 13
$receiving.Range("F16").Value = 8
$receiving.Range("G16").Value = 10
$receiving.Range("H16").Value = 7

$receiving.Range("A17").Value = 15
